$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1459.875
$ws.Range("I19").Value = 1338.1666
$ws.Range("J19").Value = 1825
$ws.Range("K19").Value = 1338.1666
$ws.Range("L19").Value = 1825
$ws.Range("M19").Value = -1163.1666
$ws.Range("N19").Value = -2175

$ws.Range("H41").Value = 692.8889
$ws.Range("I41").Value = 343
$ws.Range("J41").Value = 1392.6666
$ws.Range("K41").Value = 343
$ws.Range("L41").Value = 1392.6666
$ws.Range("M41").Value = 97
$ws.Range("N41").Value = -2272.6666

$ws.Range("H51").Value = 10462
$ws.Range("J51").Value = 10462
$ws.Range("L51").Value = 10462
$ws.Range("N51").Value = -11430

$ws.Range("H58").Value = 1700.25
$ws.Range("J58").Value = 4335.3335
$ws.Range("L58").Value = 13006.0005
$ws.Range("N58").Value = -13306.0005

$ws.Range("H74").Value = 91127.7
$ws.Range("I74").Value = 104222.31
$ws.Range("K74").Value = 104222.31
$ws.Range("M74").Value = -103286.31

$ws.Range("H77").Value = 91127.7
$ws.Range("I77").Value = 104222.31
$ws.Range("K77").Value = 521111.55
$ws.Range("M77").Value = -516431.55

$ws.Range("H80").Value = 91612.664
$ws.Range("I80").Value = 2840.5
$ws.Range("K80").Value = 8521.5
$ws.Range("M80").Value = -7523.5

$ws.Range("H83").Value = 91612.664
$ws.Range("I83").Value = 2840.5
$ws.Range("K83").Value = 25564.5
$ws.Range("M83").Value = -20572.5

$ws.Range("H92").Value = 488
$ws.Range("J92").Value = 864
$ws.Range("L92").Value = 864
$ws.Range("N92").Value = -3360

$ws.Range("H98").Value = 6333.8887
$ws.Range("I98").Value = 1002.5
$ws.Range("J98").Value = 7857.143
$ws.Range("K98").Value = 1002.5
$ws.Range("L98").Value = 7857.143
$ws.Range("M98").Value = 495.5
$ws.Range("N98").Value = -10853.143

$ws.Range("H122").Value = 6333.8887
$ws.Range("I122").Value = 1002.5
$ws.Range("J122").Value = 7857.143
$ws.Range("K122").Value = 3007.5
$ws.Range("L122").Value = 23571.429
$ws.Range("M122").Value = -557.5
$ws.Range("N122").Value = -28471.429

$ws.Range("H125").Value = 250001940
$ws.Range("I125").Value = 500001180
$ws.Range("J125").Value = 125002310
$ws.Range("K125").Value = 4500010620
$ws.Range("L125").Value = 1125020790
$ws.Range("M125").Value = -4500008160
$ws.Range("N125").Value = -1125025710

$ws.Range("H138").Value = 3785
$ws.Range("I138").Value = 2873.9333
$ws.Range("J138").Value = 4218.8413
$ws.Range("K138").Value = 8621.7999
$ws.Range("L138").Value = 12656.5239
$ws.Range("M138").Value = -3481.7999
$ws.Range("N138").Value = -22936.5239

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("M117").ClearContents()
$ws.Range("N117").ClearContents()

$ws.Range("H122").Value = 1934.125
$ws.Range("I122").Value = 1094.6
$ws.Range("K122").Value = 3283.8
$ws.Range("M122").Value = -833.7999999999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1115.8
$ws.Range("J22").Value = 1066
$ws.Range("L22").Value = 1066
$ws.Range("N22").Value = -1412

$ws.Range("H99").Value = 2983.476
$ws.Range("I99").Value = 2269.5334
$ws.Range("K99").Value = 2269.5334
$ws.Range("M99").Value = -771.5333999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H99").Value = 9246.619000000001
$ws.Range("I99").Value = 10085.615
$ws.Range("K99").Value = 10085.615
$ws.Range("M99").Value = -8587.615

$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

$ws.Range("H126").Value = 9246.619000000001
$ws.Range("I126").Value = 10085.615
$ws.Range("K126").Value = 30256.845
$ws.Range("M126").Value = -27786.845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1201.0834
$ws.Range("I46").Value = 200.28572
$ws.Range("J46").Value = 2602.2
$ws.Range("K46").Value = 600.85716
$ws.Range("L46").Value = 7806.599999999999
$ws.Range("M46").Value = -509.85716
$ws.Range("N46").Value = -7988.599999999999

$ws.Range("H64").Value = 14057
$ws.Range("I64").Value = 699.5
$ws.Range("J64").Value = 19400
$ws.Range("K64").Value = 2098.5
$ws.Range("L64").Value = 58200
$ws.Range("M64").Value = -1828.5
$ws.Range("N64").Value = -58740

$ws.Range("H67").Value = 14057
$ws.Range("I67").Value = 699.5
$ws.Range("J67").Value = 19400
$ws.Range("K67").Value = 2098.5
$ws.Range("L67").Value = 58200
$ws.Range("M67").Value = -1162.5
$ws.Range("N67").Value = -60072

$ws.Range("H132").Value = 2788.2354
$ws.Range("I132").Value = 2093.3333
$ws.Range("K132").Value = 18839.9997
$ws.Range("M132").Value = -16309.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5303
$ws.Range("I70").Value = 4872.375
$ws.Range("J70").Value = 5795.143
$ws.Range("K70").Value = 4872.375
$ws.Range("L70").Value = 5795.143
$ws.Range("M70").Value = -4602.375
$ws.Range("N70").Value = -6335.143

$ws.Range("H73").Value = 5303
$ws.Range("I73").Value = 4872.375
$ws.Range("J73").Value = 5795.143
$ws.Range("K73").Value = 4872.375
$ws.Range("L73").Value = 5795.143
$ws.Range("M73").Value = -3936.375
$ws.Range("N73").Value = -7667.143

$ws.Range("H95").Value = 45980.668
$ws.Range("J95").Value = 45980.668
$ws.Range("L95").Value = 45980.668
$ws.Range("N95").Value = -51472.668

$ws.Range("H102").Value = 2409.8823
$ws.Range("I102").Value = 2140.5715
$ws.Range("K102").Value = 2140.5715
$ws.Range("M102").Value = -518.5715

$ws.Range("H121").Value = 48895.5
$ws.Range("J121").Value = 48895.5
$ws.Range("L121").Value = 48895.5
$ws.Range("N121").Value = -52389.5

$ws.Range("H140").Value = 69992.5
$ws.Range("J140").Value = 84985
$ws.Range("L140").Value = 84985
$ws.Range("N140").Value = -95345

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 1000000
$ws.Range("I4").Value = 1000000
$ws.Range("K4").Value = 1000000
$ws.Range("M4").Value = -999887

$ws.Range("H16").Value = 1175
$ws.Range("I16").Value = 997.75
$ws.Range("K16").Value = 997.75
$ws.Range("M16").Value = -827.75

$ws.Range("H28").Value = 1000000
$ws.Range("I28").Value = 1000000
$ws.Range("K28").Value = 1000000
$ws.Range("M28").Value = -999768

$ws.Range("H37").Value = 1000000
$ws.Range("I37").Value = 1000000
$ws.Range("K37").Value = 1000000
$ws.Range("M37").Value = -999893

$ws.Range("H56").Value = 6272525.5
$ws.Range("I56").Value = 6272525.5
$ws.Range("K56").Value = 6272525.5
$ws.Range("M56").Value = -6271834.5

$ws.Range("H82").Value = 1649.6666
$ws.Range("I82").Value = 1624.75
$ws.Range("J82").Value = 1699.5
$ws.Range("K82").Value = 1624.75
$ws.Range("L82").Value = 1699.5
$ws.Range("M82").Value = -1263.75
$ws.Range("N82").Value = -2421.5

$ws.Range("H85").Value = 1649.6666
$ws.Range("I85").Value = 1624.75
$ws.Range("J85").Value = 1699.5
$ws.Range("K85").Value = 1624.75
$ws.Range("L85").Value = 1699.5
$ws.Range("M85").Value = -376.75
$ws.Range("N85").Value = -4195.5

$ws.Range("H122").Value = 7266.9653
$ws.Range("I122").Value = 7178.381
$ws.Range("J122").Value = 7499.5
$ws.Range("K122").Value = 21535.143
$ws.Range("L122").Value = 22498.5
$ws.Range("M122").Value = -19085.143
$ws.Range("N122").Value = -27398.5

$ws.Range("H130").Value = 34999
$ws.Range("J130").Value = 34999
$ws.Range("L130").Value = 34999
$ws.Range("N130").Value = -45039

$ws.Range("H133").Value = 99225
$ws.Range("I133").Value = 99225
$ws.Range("K133").Value = 99225
$ws.Range("M133").Value = -96695

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()

$ws.Range("H96").Value = 5729.143
$ws.Range("I96").Value = 5885.8
$ws.Range("J96").Value = 5642.1113
$ws.Range("K96").Value = 5885.8
$ws.Range("L96").Value = 5642.1113
$ws.Range("M96").Value = -4512.8
$ws.Range("N96").Value = -8388.1113

$ws.Range("H122").Value = 6250
$ws.Range("J122").Value = 6250
$ws.Range("L122").Value = 18750
$ws.Range("N122").Value = -23650
